$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '43.222.54'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '2.355.94'
$ws.Range("E3").Value = '  +4.60%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.00'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.648'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.91'
$ws.Range("E7").Value = '  +13.64%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  +7.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '26.96'
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.106'
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.705.90'
$ws.Range("E13").Value = '  +4.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.04'
$ws.Range("E14").Value = '  +3.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.28'
$ws.Range("E15").Value = '  +3.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.868'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '2.346.14'
$ws.Range("E17").Value = '  +3.62%  '
$ws.Range("D18").Value = '43.262.39'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000102'
$ws.Range("E19").Value = '  +4.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.32'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.25'
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '249.93'
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.75'
$ws.Range("E24").Value = '  +3.23%  '
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.98'
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.31'
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.60'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.51'
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("E31").Value = '  -4.96%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.98'
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.03'
$ws.Range("E35").Value = '  +2.39%  '
$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.56'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.43'
$ws.Range("E37").Value = '  +6.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.67'
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0254'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.91'
$ws.Range("E41").Value = '  +3.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.51'
$ws.Range("E42").Value = '  +7.93%  '
$ws.Range("E43").Value = '  +8.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.37'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.47'
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("E46").Value = '  +2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0954'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").Value = '1.444.55'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").Value = '2.576.94'
$ws.Range("E49").Value = '  +4.52%  '
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000202'
$ws.Range("E51").Value = '  -2.08%  '
